$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.399676619293871
$ws.Range("C2").Value = 0.4075566414707055
$ws.Range("D2").Value = 0.6541401328808547
$ws.Range("E2").Value = 0.2669736999609782
$ws.Range("G2").Value = 0.002461768424108221
$ws.Range("I2").Value = 0.6363252255562202
$ws.Range("J2").Value = 0.1390942464946505
$ws.Range("O2").Value = 3.729941506039637
$ws.Range("B3").Value = 1.26052660375268
$ws.Range("C3").Value = 0.3646095372368734
$ws.Range("D3").Value = 0.6438330498335176
$ws.Range("E3").Value = 0.2618987202113203
$ws.Range("G3").Value = 0.002465506629911782
$ws.Range("I3").Value = 0.6480131072965918
$ws.Range("J3").Value = 0.1357181705338419
$ws.Range("O3").Value = 3.751147857941021
$ws.Range("B4").Value = 1.175086224670792
$ws.Range("C4").Value = 0.3382288133309714
$ws.Range("D4").Value = 0.6378521179113932
$ws.Range("E4").Value = 0.2589310510746614
$ws.Range("G4").Value = 0.002467922622481407
$ws.Range("I4").Value = 0.6557530636660367
$ws.Range("J4").Value = 0.1337271187256519
$ws.Range("O4").Value = 3.767268228403168
$ws.Range("B5").Value = 1.140269798261158
$ws.Range("C5").Value = 0.3274760849561176
$ws.Range("D5").Value = 0.6355022051663468
$ws.Range("E5").Value = 0.257758965872533
$ws.Range("G5").Value = 0.00246893761199455
$ws.Range("I5").Value = 0.6590484421063927
$ws.Range("J5").Value = 0.1329362869147914
$ws.Range("O5").Value = 3.774614431604164
$ws.Range("B6").Value = 1.134488680230163
$ws.Range("C6").Value = 0.3256904716027407
$ws.Range("D6").Value = 0.6351172790505757
$ws.Range("E6").Value = 0.2575665912613374
$ws.Range("G6").Value = 0.002469107992296347
$ws.Range("I6").Value = 0.6596041578877738
$ws.Range("J6").Value = 0.1328062089088036
$ws.Range("O6").Value = 3.775881123207455
$ws.Range("B7").Value = 1.174616671011734
$ws.Range("C7").Value = 0.3380838073251766
$ws.Range("D7").Value = 0.6378200724977887
$ws.Range("E7").Value = 0.2589150931068147
$ws.Range("G7").Value = 0.002467936187712909
$ws.Range("I7").Value = 0.6557969348101835
$ws.Range("J7").Value = 0.133716370198286
$ws.Range("O7").Value = 3.76736415893518
$ws.Range("B8").Value = 1.351698986463305
$ws.Range("C8").Value = 0.3927509953437607
$ws.Range("D8").Value = 0.6505140546986183
$ws.Range("E8").Value = 0.2651930164138321
$ws.Range("G8").Value = 0.002463032360705102
$ws.Range("I8").Value = 0.6402379873128723
$ws.Range("J8").Value = 0.1379131514769867
$ws.Range("O8").Value = 3.7366086894001
$ws.Range("B9").Value = 1.698889374491898
$ws.Range("C9").Value = 0.499853644865027
$ws.Range("D9").Value = 0.6781702113357539
$ws.Range("E9").Value = 0.2786846683797251
$ws.Range("G9").Value = 0.002454369382985632
$ws.Range("I9").Value = 0.6142172968937807
$ws.Range("J9").Value = 0.1467955084937103
$ws.Range("O9").Value = 3.701002531436671
$ws.Range("B10").Value = 1.953883920975272
$ws.Range("C10").Value = 0.5784735408537358
$ws.Range("D10").Value = 0.7001836413139699
$ws.Range("E10").Value = 0.2893227543046351
$ws.Range("G10").Value = 0.002448579685808364
$ws.Range("I10").Value = 0.5978629068475456
$ws.Range("J10").Value = 0.1537240103233302
$ws.Range("O10").Value = 3.690062822983265
$ws.Range("B11").Value = 2.069861449971143
$ws.Range("C11").Value = 0.6142239390795794
$ws.Range("D11").Value = 0.7105684295805759
$ws.Range("E11").Value = 0.2943212945246501
$ws.Range("G11").Value = 0.002446069327572506
$ws.Range("I11").Value = 0.5910283331313693
$ws.Range("J11").Value = 0.1569645082764879
$ws.Range("O11").Value = 3.68842568505039
$ws.Range("B12").Value = 2.113774977440073
$ws.Range("C12").Value = 0.627759394426846
$ws.Range("D12").Value = 0.7145543272207533
$ws.Range("E12").Value = 0.2962370936087311
$ws.Range("G12").Value = 0.002445136364661621
$ws.Range("I12").Value = 0.5885277410818475
$ws.Range("J12").Value = 0.1582044268700287
$ws.Range("O12").Value = 3.688288700299296
$ws.Range("B13").Value = 2.104317646272079
$ws.Range("C13").Value = 0.6248444069995571
$ws.Range("D13").Value = 0.7136935162746738
$ws.Range("E13").Value = 0.2958234697861144
$ws.Range("G13").Value = 0.002445336511143426
$ws.Range("I13").Value = 0.5890623882384531
$ws.Range("J13").Value = 0.1579368174321019
$ws.Range("O13").Value = 3.688296678358313
$ws.Range("B14").Value = 2.07347434479874
$ws.Range("C14").Value = 0.6153375604117741
$ws.Range("D14").Value = 0.7108952809753362
$ws.Range("E14").Value = 0.2944784479332228
$ws.Range("G14").Value = 0.002445992218869653
$ws.Range("I14").Value = 0.5908208511816326
$ws.Range("J14").Value = 0.1570662600047115
$ws.Range("O14").Value = 3.688404720056269
$ws.Range("B15").Value = 2.054581278374314
$ws.Range("C15").Value = 0.6095140121768736
$ws.Range("D15").Value = 0.7091882384500821
$ws.Range("E15").Value = 0.2936575752710908
$ws.Range("G15").Value = 0.002446396155413798
$ws.Range("I15").Value = 0.5919093720213624
$ws.Range("J15").Value = 0.1565346887611128
$ws.Range("O15").Value = 3.688533877039362
$ws.Range("B16").Value = 1.946303956193674
$ws.Range("C16").Value = 0.5761368479782618
$ws.Range("D16").Value = 0.6995124392092578
$ws.Range("E16").Value = 0.288999297373401
$ws.Range("G16").Value = 0.002448746219020624
$ws.Range("I16").Value = 0.5983217813445911
$ws.Range("J16").Value = 0.1535140258216074
$ws.Range("O16").Value = 3.690237246899756
$ws.Range("B17").Value = 1.879872815707529
$ws.Range("C17").Value = 0.5556571143546876
$ws.Range("D17").Value = 0.6936716738712789
$ws.Range("E17").Value = 0.2861824206812287
$ws.Range("G17").Value = 0.002450219449664941
$ws.Range("I17").Value = 0.6024109258823209
$ws.Range("J17").Value = 0.1516837038094678
$ws.Range("O17").Value = 3.692139470671805
$ws.Range("B18").Value = 1.84166152293551
$ws.Range("C18").Value = 0.5438764014990056
$ws.Range("D18").Value = 0.6903471131984702
$ws.Range("E18").Value = 0.2845772104586501
$ws.Range("G18").Value = 0.00245107843351855
$ws.Range("I18").Value = 0.6048198206371076
$ws.Range("J18").Value = 0.150639295117827
$ws.Range("O18").Value = 3.693547732857922
$ws.Range("B19").Value = 1.828723567777274
$ws.Range("C19").Value = 0.5398874417343222
$ws.Range("D19").Value = 0.689227464375989
$ws.Range("E19").Value = 0.2840362852163665
$ws.Range("G19").Value = 0.002451371269509389
$ws.Range("I19").Value = 0.6056451931547713
$ws.Range("J19").Value = 0.1502871076017556
$ws.Range("O19").Value = 3.69407841835627
$ws.Range("B20").Value = 1.886944728419166
$ws.Range("C20").Value = 0.557837355686388
$ws.Range("D20").Value = 0.6942898213881392
$ws.Range("E20").Value = 0.2864807308942474
$ws.Range("G20").Value = 0.002450061419689365
$ws.Range("I20").Value = 0.6019697341358992
$ws.Range("J20").Value = 0.1518776807023841
$ws.Range("O20").Value = 3.691904441992079
$ws.Range("B21").Value = 2.082533909237554
$ws.Range("C21").Value = 0.618130021226591
$ws.Range("D21").Value = 0.7117157404384784
$ws.Range("E21").Value = 0.2948728897957693
$ws.Range("G21").Value = 0.002445799143237719
$ws.Range("I21").Value = 0.5903019689511027
$ws.Range("J21").Value = 0.1573216156240989
$ws.Range("O21").Value = 3.688359856018621
$ws.Range("B22").Value = 2.210334920035336
$ws.Range("C22").Value = 0.657520333809714
$ws.Range("D22").Value = 0.7234159064425398
$ws.Range("E22").Value = 0.3004915124430312
$ws.Range("G22").Value = 0.002443116366176452
$ws.Range("I22").Value = 0.5831867325581896
$ws.Range("J22").Value = 0.1609542500661547
$ws.Range("O22").Value = 3.688859636345825
$ws.Range("B23").Value = 2.142128201073547
$ws.Range("C23").Value = 0.636498434135035
$ws.Range("D23").Value = 0.7171427920265501
$ws.Range("E23").Value = 0.2974804777458715
$ws.Range("G23").Value = 0.002444538832258434
$ws.Range("I23").Value = 0.5869374089545119
$ws.Range("J23").Value = 0.159008589100651
$ws.Range("O23").Value = 3.688334280430439
$ws.Range("B24").Value = 1.883747576195276
$ws.Range("C24").Value = 0.5568516892435582
$ws.Range("D24").Value = 0.6940102529517276
$ws.Range("E24").Value = 0.2863458204676803
$ws.Range("G24").Value = 0.002450132827660592
$ws.Range("I24").Value = 0.6021690163145976
$ws.Range("J24").Value = 0.1517899592480347
$ws.Range("O24").Value = 3.69200971851572
$ws.Range("B25").Value = 1.604977188037424
$ws.Range("C25").Value = 0.4708910333923768
$ws.Range("D25").Value = 0.6703916489883568
$ws.Range("E25").Value = 0.2749078052675316
$ws.Range("G25").Value = 0.002456611524872976
$ws.Range("I25").Value = 0.6207732211878074
$ws.Range("J25").Value = 0.1443222433186975
$ws.Range("O25").Value = 3.707973653263366
